$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) contains values that look numeric (e.g. "0.3680",
# "1.000") but must stay as literal text so formatting such as trailing
# zeros is preserved exactly like the source data. Temporarily force a
# text number format on those cells before assigning the values, then put
# the style back to Normal/General so no stray formatting is left behind.
# (NumberFormat must be applied per contiguous block; multi-area/union
# ranges do not reliably keep the text coercion for later Value writes.)
$dTextRanges = @(
    $ws.Range("D2:D5"),
    $ws.Range("D7:D19"),
    $ws.Range("D21:D23"),
    $ws.Range("D25:D31"),
    $ws.Range("D33:D36"),
    $ws.Range("D39:D51")
)
foreach ($r in $dTextRanges) { $r.NumberFormat = "@" }

$ws.Range("D2").Value = "26.777.36"
$ws.Range("E2").Value = "  -0.91%  "

$ws.Range("D3").Value = "1.796.97"
$ws.Range("E3").Value = "  -1.21%  "

$ws.Range("D4").Value = "0.9998"
$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").Value = "309.26"
$ws.Range("E5").Value = "  -0.45%  "

$ws.Range("E6").Value = "  -0.02%  "

$ws.Range("D7").Value = "0.4392"
$ws.Range("E7").Value = "  +4.09%  "

$ws.Range("D8").Value = "0.3680"
$ws.Range("E8").Value = "  +0.25%  "

$ws.Range("D9").Value = "0.07386"
$ws.Range("E9").Value = "  +2.48%  "

$ws.Range("D10").Value = "0.8540"
$ws.Range("E10").Value = "  +1.80%  "

$ws.Range("D11").Value = "20.59"
$ws.Range("E11").Value = "  -0.84%  "

$ws.Range("D12").Value = "1.938.19"
$ws.Range("E12").Value = "  +6.56%  "

$ws.Range("D13").Value = "6.605"
$ws.Range("E13").Value = "  -0.58%  "

$ws.Range("D14").Value = "92.07"
$ws.Range("E14").Value = "  +3.10%  "

$ws.Range("D15").Value = "0.07047"
$ws.Range("E15").Value = "  -0.38%  "

$ws.Range("D16").Value = "5.255"
$ws.Range("E16").Value = "  -0.41%  "

$ws.Range("D17").Value = "1.001"
$ws.Range("E17").Value = "  -0.10%  "

$ws.Range("D18").Value = "0.000008653"
$ws.Range("E18").Value = "  -1.63%  "

$ws.Range("D19").Value = "1.000"
$ws.Range("E19").Value = "  -0.02%  "

$ws.Range("D21").Value = "26.807.34"
$ws.Range("E21").Value = "  -0.94%  "

$ws.Range("D22").Value = "5.142"
$ws.Range("E22").Value = "  +0.43%  "

$ws.Range("D23").Value = "10.81"
$ws.Range("E23").Value = "  -0.27%  "

$ws.Range("E24").Value = "  +0.10%  "

$ws.Range("D25").Value = "151.54"
$ws.Range("E25").Value = "  -0.16%  "

$ws.Range("D26").Value = "2.199"
$ws.Range("E26").Value = "  -0.71%  "

$ws.Range("D27").Value = "18.32"
$ws.Range("E27").Value = "  +0.50%  "

$ws.Range("D28").Value = "5.185"
$ws.Range("E28").Value = "  -0.54%  "

$ws.Range("D29").Value = "117.02"
$ws.Range("E29").Value = "  +0.69%  "

$ws.Range("D30").Value = "0.08786"
$ws.Range("E30").Value = "  +0.49%  "

$ws.Range("D31").Value = "0.7351"
$ws.Range("E31").Value = "  -0.31%  "

$ws.Range("E32").Value = "  -2.12%  "

$ws.Range("D33").Value = "2.914"
$ws.Range("E33").Value = "  -1.54%  "

$ws.Range("D34").Value = "4.429"
$ws.Range("E34").Value = "  +0.52%  "

$ws.Range("D35").Value = "0.9994"
$ws.Range("E35").Value = "  -0.07%  "

$ws.Range("D36").Value = "1.082"
$ws.Range("E36").Value = "  -0.76%  "

$ws.Range("E37").Value = "  -0.01%  "

$ws.Range("E38").Value = "  -1.15%  "

$ws.Range("D39").Value = "0.5208"
$ws.Range("E39").Value = "  +3.62%  "

$ws.Range("D40").Value = "7.005"
$ws.Range("E40").Value = "  -3.63%  "

$ws.Range("D41").Value = "2.806"
$ws.Range("E41").Value = "  -2.18%  "

$ws.Range("D42").Value = "0.1676"
$ws.Range("E42").Value = "  -0.66%  "

$ws.Range("D43").Value = "0.5006"
$ws.Range("E43").Value = "  +6.34%  "

$ws.Range("D44").Value = "8.438"
$ws.Range("E44").Value = "  -1.65%  "

$ws.Range("D45").Value = "1.960"
$ws.Range("E45").Value = "  +4.51%  "

$ws.Range("D46").Value = "10.29"
$ws.Range("E46").Value = "  -1.76%  "

$ws.Range("D47").Value = "105.01"
$ws.Range("E47").Value = "  -1.02%  "

$ws.Range("D48").Value = "0.9994"
$ws.Range("E48").Value = "  -0.01%  "

$ws.Range("D49").Value = "1.657"
$ws.Range("E49").Value = "  +0.82%  "

$ws.Range("D50").Value = "0.06313"
$ws.Range("E50").Value = "  -0.49%  "

$ws.Range("D51").Value = "0.9144"
$ws.Range("E51").Value = "  +1.65%  "

# Restore default formatting on the Price column cells we touched.
foreach ($r in $dTextRanges) {
    $r.NumberFormat = "General"
    $r.Style = "Normal"
}
